$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (B5:B13) to reflect the refactored module order ---
$ws.Range("B5").Value  = "ENG1044"
$ws.Range("B6").Value  = "NET2201"
$ws.Range("B7").Value  = "NET1014"
$ws.Range("B8").Value  = "WEB1201"
$ws.Range("B9").Value  = "PRG1203"
$ws.Range("B10").Value = "NET3204 "
$ws.Range("B11").Value = "NET2201"
$ws.Range("B12").Value = "NET1014"
$ws.Range("B13").Value = "CSC3044 "

# --- Add new rows for Semester 4 ---
$ws.Range("A14").Value = 4
$ws.Range("B14").Value = "SEG1201"
$ws.Range("A15").Value = 4
$ws.Range("B15").Value = "CSC2103"
$ws.Range("A16").Value = 4
$ws.Range("B16").Value = "WEB1201"
$ws.Range("A17").Value = 4
$ws.Range("B17").Value = "PRG1203"

# --- Add new rows for Semester 5 ---
$ws.Range("A18").Value = 5
$ws.Range("B18").Value = "CSC3024"
$ws.Range("A19").Value = 5
$ws.Range("B19").Value = "CSC3034"
$ws.Range("A20").Value = 5
$ws.Range("B20").Value = "CSC3206"
$ws.Range("A21").Value = 5
$ws.Range("B21").Value = "SEG2202"

# --- Build the "wrap text + vertical center" look on an out-of-the-way helper
#     cell first (this is the combination used for the Semester 5 rows), then
#     paste just the resulting format onto A18:B21 so only a single new style
#     entry is produced rather than one per property assignment. ---
$helper = $ws.Range("Z1")
$helper.WrapText = $true
$helper.VerticalAlignment = -4108   # xlCenter

$helper.Copy()
$target = $ws.Range("A18:B21")
$target.PasteSpecial(-4122)         # xlPasteFormats
$excel.CutCopyMode = $false

$helper.Clear()

# --- Update selection to match the saved state ---
$ws.Range("C21").Select()
